# "cable fix 'n' trix"
# The Ledningsklass (cable type) column (B) has two cable-type codes that
# were too vague ("AKKJ240" / "EKKJ10"); they get their cross-section
# suffix appended ("AKKJ240/72" / "EKKJ10/10") everywhere they occur in
# the sheet. Column widths for B/D are widened to fit, and the final
# selection is left on D109.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Whole-cell (not substring) replace so we don't clobber any other value
# that merely contains these codes as a substring.
$colB = $ws.Range("B1:B123")
$colB.Replace("AKKJ240", "AKKJ240/72", [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole) | Out-Null
$colB.Replace("EKKJ10", "EKKJ10/10", [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole) | Out-Null

# Widen column B ("Ledningsklass") and D ("Slutknut") to fit the longer
# cable-type text.
$ws.Columns.Item(2).ColumnWidth = 23.67
$ws.Columns.Item(4).ColumnWidth = 15.67

# Leave the selection where the editor ended up.
$ws.Range("D109").Select()
